$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove Logan Lightning (row 9) and Capalaba (row 4) - teams no longer in the league
$ws.Rows(9).Delete()
$ws.Rows(4).Delete()

# Insert two new rows before "SC Wanderers" (currently row 11) for the new teams
$ws.Rows(11).Insert()
$ws.Rows(11).Insert()

# Populate the new rows (order chosen to match shared-string append order)
$ws.Range("B11").Value = "Redlands United"
$ws.Range("C11").Value = "RED"
$ws.Range("C12").Value = "ROC"
$ws.Range("B12").Value = "Rochedale Rovers"
$ws.Range("A11").Value = "aarq4u59zcgcm3km6yb0rikjt"
$ws.Range("A12").Value = "4yy1apqvpdzlrugwaamnw8yoz"

# Apply the special (green) font to the new contestantId cells
$ws.Range("A11").Font.Color = 6999710
$ws.Range("A11").Font.Name = "Consolas"
$ws.Range("A11").Font.Size = 11

$ws.Range("A12").Font.Color = 6999710
$ws.Range("A12").Font.Name = "Consolas"
$ws.Range("A12").Font.Size = 11

# Update active cell selection
$ws.Range("A12").Select() | Out-Null
